# Apply the "Updated cryptos list" refresh: new Price (D) and Volume(1h) (E)
# values scraped for this run. Values are plain text cells in the sheet
# (coinranking.com price strings use "." as a thousands separator, and the
# percentage column keeps its original padding), so some Price values that
# look like plain decimals need to be pinned to Text format first so Excel
# does not silently convert them (and drop trailing zeros) to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise parse as a plain number in Excel;
# force Text number format before writing so the value stays a string.
$textCells = @("D5", "D6", "D7", "D9", "D14", "D16", "D20", "D22", "D25", "D29", "D36", "D38", "D40", "D43", "D47")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.826.29'
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("D3").Value = '2.353.79'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '0.671'
$ws.Range("E5").Value = '  +3.56%  '
$ws.Range("D6").Value = '237.12'
$ws.Range("E6").Value = '  +2.09%  '
$ws.Range("D7").Value = '72.97'
$ws.Range("E7").Value = '  +10.85%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '0.540'
$ws.Range("E9").Value = '  +19.48%  '
$ws.Range("E10").Value = '  +3.01%  '
$ws.Range("E11").Value = '  +6.16%  '
$ws.Range("E12").Value = '  +2.44%  '
$ws.Range("D13").Value = '2.699.46'
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").Value = '16.67'
$ws.Range("E14").Value = '  +8.25%  '
$ws.Range("E15").Value = '  +6.16%  '
$ws.Range("D16").Value = '0.895'
$ws.Range("E16").Value = '  +5.64%  '
$ws.Range("D17").Value = '2.359.54'
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").Value = '43.820.01'
$ws.Range("E19").Value = '  +3.06%  '
$ws.Range("D20").Value = '77.89'
$ws.Range("E20").Value = '  +5.14%  '
$ws.Range("E21").Value = '  +3.43%  '
$ws.Range("D22").Value = '253.86'
$ws.Range("E22").Value = '  +1.85%  '
$ws.Range("E24").Value = '  -2.52%  '
$ws.Range("D25").Value = '2.50'
$ws.Range("E25").Value = '  +3.00%  '
$ws.Range("E26").Value = '  +6.29%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("D29").Value = '172.61'
$ws.Range("E30").Value = '  +6.09%  '
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("E32").Value = '  +4.57%  '
$ws.Range("E33").Value = '  +3.61%  '
$ws.Range("E34").Value = '  +3.61%  '
$ws.Range("E35").Value = '  +4.33%  '
$ws.Range("D36").Value = '4.05'
$ws.Range("E36").Value = '  +12.08%  '
$ws.Range("E37").Value = '  -4.55%  '
$ws.Range("D38").Value = '6.39'
$ws.Range("E38").Value = '  -0.83%  '
$ws.Range("E39").Value = '  +6.34%  '
$ws.Range("D40").Value = '19.66'
$ws.Range("E40").Value = '  +8.72%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("E42").Value = '  -2.18%  '
$ws.Range("D43").Value = '1.24'
$ws.Range("E43").Value = '  +3.47%  '
$ws.Range("E44").Value = '  +3.46%  '
$ws.Range("E45").Value = '  -0.92%  '
$ws.Range("E46").Value = '  +1.25%  '
$ws.Range("D47").Value = '97.74'
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("E48").Value = '  +11.47%  '
$ws.Range("E49").Value = '  +2.48%  '
$ws.Range("D50").Value = '1.434.35'
$ws.Range("E50").Value = '  -0.22%  '
$ws.Range("E51").Value = '  +1.53%  '
